$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '37.115.28'
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.20%  '
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.050.51'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  -1.21%  '
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.04%  '
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '248.61'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  -1.97%  '
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -1.84%  '
$c.Style = 'Normal'

$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '58.42'
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  -6.09%  '
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  -0.31%  '
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  -3.05%  '
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +0.08%  '
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '15.72'
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  -0.48%  '
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '2.350.97'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  -1.18%  '
$c.Style = 'Normal'

$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '0.831'
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +0.66%  '
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '5.76'
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +5.17%  '
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '2.052.36'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -1.21%  '
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '17.80'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +17.42%  '
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '37.128.20'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  -0.16%  '
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '74.84'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  +0.00%  '
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '0.0₃0899'
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -2.61%  '
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  -2.37%  '
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '237.29'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -1.26%  '
$c.Style = 'Normal'

$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '2.47'
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +2.00%  '
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.19'
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -1.71%  '
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '169.17'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  -1.36%  '
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.36'
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +0.90%  '
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  -1.54%  '
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  -1.41%  '
$c.Style = 'Normal'

$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.81'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +0.68%  '
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +2.26%  '
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  -2.52%  '
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '4.55'
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +2.77%  '
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0904'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +0.43%  '
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  -0.18%  '
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  -1.90%  '
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  -0.63%  '
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  -2.03%  '
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  -4.00%  '
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +10.89%  '
$c.Style = 'Normal'

$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '5.14'
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +10.25%  '
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -2.05%  '
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '17.39'
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  -5.19%  '
$c.Style = 'Normal'

$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '1.14'
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -1.59%  '
$c.Style = 'Normal'

$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '96.05'
$c.Style = 'Normal'

$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  -2.58%  '
$c.Style = 'Normal'

$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.46'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  -2.09%  '
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '1.285.90'
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  -1.96%  '
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  -1.68%  '
$c.Style = 'Normal'

$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '6.85'
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  -1.38%  '
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '2.243.28'
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  -0.68%  '
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '3.53'
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -20.68%  '
$c.Style = 'Normal'
